$wb = $excel.ActiveWorkbook

# --- parameters sheet: add "description" column values (col C) for
#     sigma, delta, gamma rows, and min/max (cols D/E) for the tau row ---
$wsParams = $wb.Worksheets.Item("parameters")

$wsParams.Range("C3").Value = "rate of becoming symptomatic per unit time"
$wsParams.Range("C2").Value = "rate of transition (exposed to pre-symptomatic)"
$wsParams.Range("C4").Value = "Recovery rate per unit time"
$wsParams.Range("C5").Value = "reduction in contact rate"
$wsParams.Range("D5").Value = 0
$wsParams.Range("E5").Value = 1

# --- update the last selected/active cell on a few sheets to match
#     where the author's cursor ended up when they saved the file ---
$wsParams.Activate()
$wsParams.Range("F12").Select() | Out-Null

$wsJurisdiction = $wb.Worksheets.Item("jurisdiction")
$wsJurisdiction.Activate()
$wsJurisdiction.Range("B9").Select() | Out-Null

$wsTravel = $wb.Worksheets.Item("travel")
$wsTravel.Activate()
$wsTravel.Range("A5").Select() | Out-Null

# restore original active sheet (parameters, second tab) as the tab shown
$wsParams.Activate()
